$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the x location value; dependent formulas in column D recalc automatically.
$ws.Range("C2").Value = 75

$excel.CalculateFullRebuild()
